# Daily attendance processing - 2025-10-18 15:39:50
# Re-order the "Recorded By" (column G) value so that "System" is listed
# last instead of first, e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System" (same for the admin@admin.com account).
# Rows that combine System with more than one other account (e.g.
# "System, system, backup@backdoor.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
        $changed++
    } elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
        $changed++
    }
}

Write-Host "Recorded By cells reordered:" $changed
